# Update Handback status report timestamps and the "ht" -> "mt" status value
# as described by the commit "Generate Report for Handback".

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2017-02-21 11:26:20"
$wsOverview.Range("G3").Value = "2017-02-21 11:26:20"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H2").Value = "2017-02-21 11:26:03"
$wsZhCn.Range("H3").Value = "2017-02-21 11:26:03"
$wsZhCn.Range("L2").Value = "2017-02-21 11:27:01"
$wsZhCn.Range("L3").Value = "2017-02-21 11:27:01"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H2").Value = "2017-02-21 11:26:20"
$wsDeDe.Range("H3").Value = "2017-02-21 11:26:20"
$wsDeDe.Range("L2").Value = "2017-02-21 11:27:24"
$wsDeDe.Range("L3").Value = "2017-02-21 11:27:24"
